$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string shows up on every sheet (Overview E/F columns, and the
#    "Status" column (C) on the zh-cn / de-de detail sheets).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in "Latest Target File" (I) / "Latest Handback File" (J)
#    for both rows, re-point the handback datetime (K) to the real value, and
#    add hyperlinks on column I matching the ones already on column A.
# ---------------------------------------------------------------------------
$wsZhCn.Range("I2").Value = "46c59b38-1b92-40fb-b875-cd8402a6dc5a.md"
$wsZhCn.Range("J2").Value = "46c59b38-1b92-40fb-b875-cd8402a6dc5a.4b590258a687be1921753c21a081e4ac8c5105fb.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-10-27 10:10:50"

$wsZhCn.Range("I3").Value = "f94e23ab-021f-4b43-94c2-c89c08623796.md"
$wsZhCn.Range("J3").Value = "f94e23ab-021f-4b43-94c2-c89c08623796.4887b6ab5826227363bd668c041ba686e9678d8f.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-10-27 10:10:50"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/428bcd75e31aa4ee452ca1d16fb0a1bc9312f4cc/e2e/46c59b38-1b92-40fb-b875-cd8402a6dc5a.md", [Type]::Missing, [Type]::Missing, "46c59b38-1b92-40fb-b875-cd8402a6dc5a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/428bcd75e31aa4ee452ca1d16fb0a1bc9312f4cc/e2e/46c59b38-1b92-40fb-b875-cd8402a6dc5a.md", [Type]::Missing, [Type]::Missing, "46c59b38-1b92-40fb-b875-cd8402a6dc5a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/428bcd75e31aa4ee452ca1d16fb0a1bc9312f4cc/e2e/f94e23ab-021f-4b43-94c2-c89c08623796.md", [Type]::Missing, [Type]::Missing, "f94e23ab-021f-4b43-94c2-c89c08623796.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/428bcd75e31aa4ee452ca1d16fb0a1bc9312f4cc/e2e/f94e23ab-021f-4b43-94c2-c89c08623796.md", [Type]::Missing, [Type]::Missing, "f94e23ab-021f-4b43-94c2-c89c08623796.md")

# Match the existing hyperlink text style on the newly-linked cells.
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 13071615
$wsZhCn.Range("I3").Font.Underline = 2
$wsZhCn.Range("I3").Font.Color = 13071615

# Column widths widened to fit the new handback-file text.
$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777050018311
$wsZhCn.Columns.Item(9).ColumnWidth = 40
$wsZhCn.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# 3. de-de sheet: same shape of update, but this handback finished later so
#    the datetime + target xlf names are their own de-de specific values.
# ---------------------------------------------------------------------------
$wsDeDe.Range("I2").Value = "46c59b38-1b92-40fb-b875-cd8402a6dc5a.md"
$wsDeDe.Range("J2").Value = "46c59b38-1b92-40fb-b875-cd8402a6dc5a.4b590258a687be1921753c21a081e4ac8c5105fb.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-10-27 10:11:09"

$wsDeDe.Range("I3").Value = "f94e23ab-021f-4b43-94c2-c89c08623796.md"
$wsDeDe.Range("J3").Value = "f94e23ab-021f-4b43-94c2-c89c08623796.4887b6ab5826227363bd668c041ba686e9678d8f.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-10-27 10:11:09"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/428bcd75e31aa4ee452ca1d16fb0a1bc9312f4cc/e2e/46c59b38-1b92-40fb-b875-cd8402a6dc5a.md", [Type]::Missing, [Type]::Missing, "46c59b38-1b92-40fb-b875-cd8402a6dc5a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/428bcd75e31aa4ee452ca1d16fb0a1bc9312f4cc/e2e/46c59b38-1b92-40fb-b875-cd8402a6dc5a.md", [Type]::Missing, [Type]::Missing, "46c59b38-1b92-40fb-b875-cd8402a6dc5a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/428bcd75e31aa4ee452ca1d16fb0a1bc9312f4cc/e2e/f94e23ab-021f-4b43-94c2-c89c08623796.md", [Type]::Missing, [Type]::Missing, "f94e23ab-021f-4b43-94c2-c89c08623796.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/428bcd75e31aa4ee452ca1d16fb0a1bc9312f4cc/e2e/f94e23ab-021f-4b43-94c2-c89c08623796.md", [Type]::Missing, [Type]::Missing, "f94e23ab-021f-4b43-94c2-c89c08623796.md")

$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = 13071615
$wsDeDe.Range("I3").Font.Underline = 2
$wsDeDe.Range("I3").Font.Color = 13071615

$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777050018311
$wsDeDe.Columns.Item(9).ColumnWidth = 40
$wsDeDe.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# 4. Overview sheet column widths (zh-cn / de-de columns widened to match).
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777050018311
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777050018311

Write-Output "Handback report generated"
